# Mandatory Field validation run sucessful
# - Widen column C (from ~23.57 to ~55.71 characters)
# - Update the active selection/view to cell C4 (also clears the scrolled
#   "topLeftCell" position that was left over from A26)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column C. The host's ColumnWidth setter quantizes to the nearest
# 1/6 of a character, so feed it the value whose round-trip lands closest
# to the desired stored width of 55.7109375 characters.
$ws.Columns.Item(3).ColumnWidth = 54.877604166666664

# Move the selection/active cell to C4 and scroll it into view; this also
# drops the stale topLeftCell="A8" left over from the previous selection.
$ws.Range("C4").Select()
